$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update test result values
$ws.Range("E3").Value = "PASS"
$ws.Range("E25").Value = "SKIP"

# Make this sheet active and update the view / selection state
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D2:D25").Select()
